$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (matches original inlineStr text cells) then set the new values.
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.999999715665001"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.879140052591499"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.995088226113807"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.994500027996446"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.999766160128656"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.999724347062487"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.867173527352789"
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = "0.999449251287372"
$ws.Range("J14").NumberFormat = "@"
$ws.Range("J14").Value = "0.000258409368966036"
$ws.Range("L14").NumberFormat = "@"
$ws.Range("L14").Value = "1228.89159505488"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.97705864973909"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.936265166036628"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.999687578373293"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.994676205886985"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.998331913118766"
$ws.Range("F23").NumberFormat = "@"
$ws.Range("F23").Value = "0.995242587915354"
$ws.Range("J23").NumberFormat = "@"
$ws.Range("J23").Value = "0.000307410279866116"
$ws.Range("L23").NumberFormat = "@"
$ws.Range("L23").Value = "2708.25997279809"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.998067649401291"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.920346057176682"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "0.889412740332115"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "0.958427942433311"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "0.927842354321809"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "0.90451852803558"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "0.987023914145289"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "0.969797598368199"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "0.999327048148483"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.921502233691681"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.94920801306659"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.977856249753898"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.987107167321343"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "0.974689176594734"
$ws.Range("F44").NumberFormat = "@"
$ws.Range("F44").Value = "0.999735693644705"
$ws.Range("J44").NumberFormat = "@"
$ws.Range("J44").Value = "0.000138994748652222"
$ws.Range("L44").NumberFormat = "@"
$ws.Range("L44").Value = "1479.52682213361"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.898996779633501"
$ws.Range("F48").NumberFormat = "@"
$ws.Range("F48").Value = "0.999641647347348"
$ws.Range("J48").NumberFormat = "@"
$ws.Range("J48").Value = "0.000245417552853434"
$ws.Range("L48").NumberFormat = "@"
$ws.Range("L48").Value = "2647.68236060822"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.945430671093185"
